# psidis/expdata/20015.xlsx -- "Add files via upload"
#
# Re-upload of the dataset sheet: the "target" column value used for
# these HERMES rows was shortened from "helium" to "h", and the header
# row (row 1) was given a bold + centered style. The previously-selected
# cell is also updated to reflect where the author's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (target) for the data rows: "helium" -> "h"
$targetRange = $ws.Range("G2:G10")
for ($i = 1; $i -le $targetRange.Cells.Count; $i++) {
    $cell = $targetRange.Cells.Item($i)
    if ($cell.Value2 -eq "helium") {
        $cell.Value = "h"
    }
}

# Header row (A1:K1): bold + centered
$headerRange = $ws.Range("A1:K1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  ## xlCenter

# Move the active selection to K16, matching the saved view state
[void]$ws.Range("K16").Select()
